# Update cryptos price/volume data as of Wed Aug 16 13:24:24 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains its original text representation
# instead of being auto-converted to a number by Excel when assigned.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.236.92'
$ws.Range("E2").Value = '  -0.47%  '

$ws.Range("D3").Value = '1.828.30'
$ws.Range("E3").Value = '  -0.78%  '

$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.43%  '

$ws.Range("D5").Value = '234.62'
$ws.Range("E5").Value = '  -1.91%  '

$ws.Range("D6").Value = '0.5979'
$ws.Range("E6").Value = '  -4.62%  '

$ws.Range("D7").Value = '1.005'
$ws.Range("E7").Value = '  +0.44%  '

$ws.Range("D8").Value = '0.06970'
$ws.Range("E8").Value = '  -5.96%  '

$ws.Range("D9").Value = '0.2755'
$ws.Range("E9").Value = '  -4.81%  '

$ws.Range("D10").Value = '23.30'
$ws.Range("E10").Value = '  -6.24%  '

$ws.Range("D11").Value = '0.07632'
$ws.Range("E11").Value = '  -1.05%  '

$ws.Range("D12").Value = '1.837.88'
$ws.Range("E12").Value = '  -0.23%  '

$ws.Range("D13").Value = '4.767'
$ws.Range("E13").Value = '  -4.04%  '

$ws.Range("D14").Value = '0.6258'
$ws.Range("E14").Value = '  -7.31%  '

$ws.Range("D15").Value = '0.000009761'
$ws.Range("E15").Value = '  -4.77%  '

$ws.Range("D16").Value = '78.53'
$ws.Range("E16").Value = '  -4.08%  '

$ws.Range("D17").Value = '29.036.72'
$ws.Range("E17").Value = '  -1.08%  '

$ws.Range("D18").Value = '5.736'
$ws.Range("E18").Value = '  -8.42%  '

$ws.Range("D19").Value = '222.40'
$ws.Range("E19").Value = '  -5.07%  '

$ws.Range("E20").Value = '  +0.38%  '

$ws.Range("D21").Value = '11.57'
$ws.Range("E21").Value = '  -6.10%  '

$ws.Range("D22").Value = '6.889'
$ws.Range("E22").Value = '  -5.85%  '

$ws.Range("D23").Value = '1.005'
$ws.Range("E23").Value = '  +0.47%  '

$ws.Range("D24").Value = '156.17'
$ws.Range("E24").Value = '  -1.05%  '

$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").Value = '7.957'
$ws.Range("E25").Value = '  -6.23%  '

$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").Value = '0.1291'
$ws.Range("E26").Value = '  -4.10%  '

$ws.Range("D27").Value = '16.49'
$ws.Range("E27").Value = '  -4.74%  '

$ws.Range("D28").Value = '0.06730'
$ws.Range("E28").Value = '  -7.86%  '

$ws.Range("D29").Value = '1.446'
$ws.Range("E29").Value = '  -2.09%  '

$ws.Range("D30").Value = '1.440'
$ws.Range("E30").Value = '  -2.58%  '

$ws.Range("D31").Value = '3.831'
$ws.Range("E31").Value = '  -4.92%  '

$ws.Range("D32").Value = '3.761'
$ws.Range("E32").Value = '  -7.06%  '

$ws.Range("D33").Value = '1.092'
$ws.Range("E33").Value = '  -4.55%  '

$ws.Range("D34").Value = '1.718'
$ws.Range("E34").Value = '  -5.51%  '

$ws.Range("D35").Value = '0.6428'
$ws.Range("E35").Value = '  -8.09%  '

$ws.Range("D36").Value = '2.548'
$ws.Range("E36").Value = '  -0.90%  '

$ws.Range("D37").Value = '2.730'
$ws.Range("E37").Value = '  -2.32%  '

$ws.Range("D38").Value = '1.191.93'
$ws.Range("E38").Value = '  -3.38%  '

$ws.Range("D39").Value = '0.01737'
$ws.Range("E39").Value = '  -5.34%  '

$ws.Range("D40").Value = '6.502'
$ws.Range("E40").Value = '  -5.98%  '

$ws.Range("D41").Value = '0.9019'
$ws.Range("E41").Value = '  -4.72%  '

$ws.Range("D42").Value = '1.005'
$ws.Range("E42").Value = '  +0.45%  '

$ws.Range("D43").Value = '1.983.00'
$ws.Range("E43").Value = '  -0.37%  '

$ws.Range("D44").Value = '100.38'
$ws.Range("E44").Value = '  -0.54%  '

$ws.Range("D45").Value = '62.01'
$ws.Range("E45").Value = '  -4.94%  '

$ws.Range("E46").Value = '  -4.86%  '

$ws.Range("D47").Value = '8.481'
$ws.Range("E47").Value = '  -4.09%  '

$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").Value = '0.4554'
$ws.Range("E48").Value = '  -0.37%  '

$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '0.05507'
$ws.Range("E49").Value = '  -2.66%  '

$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = '1.572'
$ws.Range("E50").Value = '  -8.06%  '

$ws.Range("D51").Value = '6.363'
$ws.Range("E51").Value = '  -8.57%  '
